# Thay đổi giao diện thời khóa biểu:
# - Xóa sheet "Tổng quan"
# - Đổi tên sheet "Chi tiết lớp học" thành "Phương án 1"

$wb = $excel.ActiveWorkbook

$excel.DisplayAlerts = $false

$wsOverview = $wb.Worksheets.Item("Tổng quan")
[void]$wsOverview.Delete()

$wsDetail = $wb.Worksheets.Item("Chi tiết lớp học")
$wsDetail.Name = "Phương án 1"

$excel.DisplayAlerts = $true
